$wb = $excel.ActiveWorkbook

# --- mmWave(InBed): append rows 100-115 ---
$wsInBed = $wb.Worksheets.Item("mmWave(InBed)")
$wsInBed.Cells.Item(100, 1).NumberFormat = "@"
$wsInBed.Cells.Item(100, 1).Value = "2026-02-01"
$wsInBed.Cells.Item(100, 2).Value = "21:06:42"
$wsInBed.Cells.Item(100, 3).Value = "21:00"
$wsInBed.Cells.Item(100, 4).Value = "Bedroom"
$wsInBed.Cells.Item(100, 5).Value = "In Bed"
$wsInBed.Cells.Item(100, 6).Value = "Occupied"

$wsInBed.Cells.Item(101, 1).NumberFormat = "@"
$wsInBed.Cells.Item(101, 1).Value = "2026-02-01"
$wsInBed.Cells.Item(101, 2).Value = "21:06:47"
$wsInBed.Cells.Item(101, 3).Value = "21:00"
$wsInBed.Cells.Item(101, 4).Value = "Bedroom"
$wsInBed.Cells.Item(101, 5).Value = "In Bed"
$wsInBed.Cells.Item(101, 6).Value = "Occupied"

$wsInBed.Cells.Item(102, 1).NumberFormat = "@"
$wsInBed.Cells.Item(102, 1).Value = "2026-02-01"
$wsInBed.Cells.Item(102, 2).Value = "21:06:54"
$wsInBed.Cells.Item(102, 3).Value = "21:00"
$wsInBed.Cells.Item(102, 4).Value = "Bedroom"
$wsInBed.Cells.Item(102, 5).Value = "In Bed"
$wsInBed.Cells.Item(102, 6).Value = "Occupied"

$wsInBed.Cells.Item(103, 1).NumberFormat = "@"
$wsInBed.Cells.Item(103, 1).Value = "2026-02-01"
$wsInBed.Cells.Item(103, 2).Value = "21:06:59"
$wsInBed.Cells.Item(103, 3).Value = "21:00"
$wsInBed.Cells.Item(103, 4).Value = "Bedroom"
$wsInBed.Cells.Item(103, 5).Value = "In Bed"
$wsInBed.Cells.Item(103, 6).Value = "Occupied"

$wsInBed.Cells.Item(104, 1).NumberFormat = "@"
$wsInBed.Cells.Item(104, 1).Value = "2026-02-01"
$wsInBed.Cells.Item(104, 2).Value = "21:07:00"
$wsInBed.Cells.Item(104, 3).Value = "21:00"
$wsInBed.Cells.Item(104, 4).Value = "Bedroom"
$wsInBed.Cells.Item(104, 5).Value = "In Bed"
$wsInBed.Cells.Item(104, 6).Value = "Occupied"

$wsInBed.Cells.Item(105, 1).NumberFormat = "@"
$wsInBed.Cells.Item(105, 1).Value = "2026-02-01"
$wsInBed.Cells.Item(105, 2).Value = "21:07:03"
$wsInBed.Cells.Item(105, 3).Value = "21:00"
$wsInBed.Cells.Item(105, 4).Value = "Bedroom"
$wsInBed.Cells.Item(105, 5).Value = "In Bed"
$wsInBed.Cells.Item(105, 6).Value = "Occupied"

$wsInBed.Cells.Item(106, 1).NumberFormat = "@"
$wsInBed.Cells.Item(106, 1).Value = "2026-02-01"
$wsInBed.Cells.Item(106, 2).Value = "21:07:05"
$wsInBed.Cells.Item(106, 3).Value = "21:00"
$wsInBed.Cells.Item(106, 4).Value = "Bedroom"
$wsInBed.Cells.Item(106, 5).Value = "In Bed"
$wsInBed.Cells.Item(106, 6).Value = "Occupied"

$wsInBed.Cells.Item(107, 1).NumberFormat = "@"
$wsInBed.Cells.Item(107, 1).Value = "2026-02-01"
$wsInBed.Cells.Item(107, 2).Value = "21:07:08"
$wsInBed.Cells.Item(107, 3).Value = "21:00"
$wsInBed.Cells.Item(107, 4).Value = "Bedroom"
$wsInBed.Cells.Item(107, 5).Value = "In Bed"
$wsInBed.Cells.Item(107, 6).Value = "Occupied"

$wsInBed.Cells.Item(108, 1).NumberFormat = "@"
$wsInBed.Cells.Item(108, 1).Value = "2026-02-01"
$wsInBed.Cells.Item(108, 2).Value = "21:07:12"
$wsInBed.Cells.Item(108, 3).Value = "21:00"
$wsInBed.Cells.Item(108, 4).Value = "Bedroom"
$wsInBed.Cells.Item(108, 5).Value = "In Bed"
$wsInBed.Cells.Item(108, 6).Value = "Occupied"

$wsInBed.Cells.Item(109, 1).NumberFormat = "@"
$wsInBed.Cells.Item(109, 1).Value = "2026-02-01"
$wsInBed.Cells.Item(109, 2).Value = "21:07:13"
$wsInBed.Cells.Item(109, 3).Value = "21:00"
$wsInBed.Cells.Item(109, 4).Value = "Bedroom"
$wsInBed.Cells.Item(109, 5).Value = "In Bed"
$wsInBed.Cells.Item(109, 6).Value = "Occupied"

$wsInBed.Cells.Item(110, 1).NumberFormat = "@"
$wsInBed.Cells.Item(110, 1).Value = "2026-02-01"
$wsInBed.Cells.Item(110, 2).Value = "21:07:14"
$wsInBed.Cells.Item(110, 3).Value = "21:00"
$wsInBed.Cells.Item(110, 4).Value = "Bedroom"
$wsInBed.Cells.Item(110, 5).Value = "In Bed"
$wsInBed.Cells.Item(110, 6).Value = "Occupied"

$wsInBed.Cells.Item(111, 1).NumberFormat = "@"
$wsInBed.Cells.Item(111, 1).Value = "2026-02-01"
$wsInBed.Cells.Item(111, 2).Value = "21:07:20"
$wsInBed.Cells.Item(111, 3).Value = "21:00"
$wsInBed.Cells.Item(111, 4).Value = "Bedroom"
$wsInBed.Cells.Item(111, 5).Value = "In Bed"
$wsInBed.Cells.Item(111, 6).Value = "Occupied"

$wsInBed.Cells.Item(112, 1).NumberFormat = "@"
$wsInBed.Cells.Item(112, 1).Value = "2026-02-01"
$wsInBed.Cells.Item(112, 2).Value = "21:07:27"
$wsInBed.Cells.Item(112, 3).Value = "21:00"
$wsInBed.Cells.Item(112, 4).Value = "Bedroom"
$wsInBed.Cells.Item(112, 5).Value = "In Bed"
$wsInBed.Cells.Item(112, 6).Value = "Occupied"

$wsInBed.Cells.Item(113, 1).NumberFormat = "@"
$wsInBed.Cells.Item(113, 1).Value = "2026-02-01"
$wsInBed.Cells.Item(113, 2).Value = "21:07:30"
$wsInBed.Cells.Item(113, 3).Value = "21:00"
$wsInBed.Cells.Item(113, 4).Value = "Bedroom"
$wsInBed.Cells.Item(113, 5).Value = "In Bed"
$wsInBed.Cells.Item(113, 6).Value = "Occupied"

$wsInBed.Cells.Item(114, 1).NumberFormat = "@"
$wsInBed.Cells.Item(114, 1).Value = "2026-02-01"
$wsInBed.Cells.Item(114, 2).Value = "21:07:36"
$wsInBed.Cells.Item(114, 3).Value = "21:00"
$wsInBed.Cells.Item(114, 4).Value = "Bedroom"
$wsInBed.Cells.Item(114, 5).Value = "In Bed"
$wsInBed.Cells.Item(114, 6).Value = "Occupied"

$wsInBed.Cells.Item(115, 1).NumberFormat = "@"
$wsInBed.Cells.Item(115, 1).Value = "2026-02-01"
$wsInBed.Cells.Item(115, 2).Value = "21:07:40"
$wsInBed.Cells.Item(115, 3).Value = "21:00"
$wsInBed.Cells.Item(115, 4).Value = "Bedroom"
$wsInBed.Cells.Item(115, 5).Value = "In Bed"
$wsInBed.Cells.Item(115, 6).Value = "Occupied"

# --- mmWave(BR): append rows 96-111 ---
$wsBR = $wb.Worksheets.Item("mmWave(BR)")
$wsBR.Cells.Item(96, 1).NumberFormat = "@"
$wsBR.Cells.Item(96, 1).Value = "2026-02-01"
$wsBR.Cells.Item(96, 2).Value = "21:06:43"
$wsBR.Cells.Item(96, 3).Value = "21:00"
$wsBR.Cells.Item(96, 4).Value = "Bedroom"
$wsBR.Cells.Item(96, 5).Value = 2
$wsBR.Cells.Item(96, 6).Value = "Occupied"

$wsBR.Cells.Item(97, 1).NumberFormat = "@"
$wsBR.Cells.Item(97, 1).Value = "2026-02-01"
$wsBR.Cells.Item(97, 2).Value = "21:06:48"
$wsBR.Cells.Item(97, 3).Value = "21:00"
$wsBR.Cells.Item(97, 4).Value = "Bedroom"
$wsBR.Cells.Item(97, 5).Value = 1
$wsBR.Cells.Item(97, 6).Value = "Occupied"

$wsBR.Cells.Item(98, 1).NumberFormat = "@"
$wsBR.Cells.Item(98, 1).Value = "2026-02-01"
$wsBR.Cells.Item(98, 2).Value = "21:06:55"
$wsBR.Cells.Item(98, 3).Value = "21:00"
$wsBR.Cells.Item(98, 4).Value = "Bedroom"
$wsBR.Cells.Item(98, 5).Value = 2
$wsBR.Cells.Item(98, 6).Value = "Occupied"

$wsBR.Cells.Item(99, 1).NumberFormat = "@"
$wsBR.Cells.Item(99, 1).Value = "2026-02-01"
$wsBR.Cells.Item(99, 2).Value = "21:07:00"
$wsBR.Cells.Item(99, 3).Value = "21:00"
$wsBR.Cells.Item(99, 4).Value = "Bedroom"
$wsBR.Cells.Item(99, 5).Value = 3
$wsBR.Cells.Item(99, 6).Value = "Occupied"

$wsBR.Cells.Item(100, 1).NumberFormat = "@"
$wsBR.Cells.Item(100, 1).Value = "2026-02-01"
$wsBR.Cells.Item(100, 2).Value = "21:07:01"
$wsBR.Cells.Item(100, 3).Value = "21:00"
$wsBR.Cells.Item(100, 4).Value = "Bedroom"
$wsBR.Cells.Item(100, 5).Value = 2
$wsBR.Cells.Item(100, 6).Value = "Occupied"

$wsBR.Cells.Item(101, 1).NumberFormat = "@"
$wsBR.Cells.Item(101, 1).Value = "2026-02-01"
$wsBR.Cells.Item(101, 2).Value = "21:07:04"
$wsBR.Cells.Item(101, 3).Value = "21:00"
$wsBR.Cells.Item(101, 4).Value = "Bedroom"
$wsBR.Cells.Item(101, 5).Value = 1
$wsBR.Cells.Item(101, 6).Value = "Occupied"

$wsBR.Cells.Item(102, 1).NumberFormat = "@"
$wsBR.Cells.Item(102, 1).Value = "2026-02-01"
$wsBR.Cells.Item(102, 2).Value = "21:07:06"
$wsBR.Cells.Item(102, 3).Value = "21:00"
$wsBR.Cells.Item(102, 4).Value = "Bedroom"
$wsBR.Cells.Item(102, 5).Value = 2
$wsBR.Cells.Item(102, 6).Value = "Occupied"

$wsBR.Cells.Item(103, 1).NumberFormat = "@"
$wsBR.Cells.Item(103, 1).Value = "2026-02-01"
$wsBR.Cells.Item(103, 2).Value = "21:07:09"
$wsBR.Cells.Item(103, 3).Value = "21:00"
$wsBR.Cells.Item(103, 4).Value = "Bedroom"
$wsBR.Cells.Item(103, 5).Value = 1
$wsBR.Cells.Item(103, 6).Value = "Occupied"

$wsBR.Cells.Item(104, 1).NumberFormat = "@"
$wsBR.Cells.Item(104, 1).Value = "2026-02-01"
$wsBR.Cells.Item(104, 2).Value = "21:07:13"
$wsBR.Cells.Item(104, 3).Value = "21:00"
$wsBR.Cells.Item(104, 4).Value = "Bedroom"
$wsBR.Cells.Item(104, 5).Value = 2
$wsBR.Cells.Item(104, 6).Value = "Occupied"

$wsBR.Cells.Item(105, 1).NumberFormat = "@"
$wsBR.Cells.Item(105, 1).Value = "2026-02-01"
$wsBR.Cells.Item(105, 2).Value = "21:07:14"
$wsBR.Cells.Item(105, 3).Value = "21:00"
$wsBR.Cells.Item(105, 4).Value = "Bedroom"
$wsBR.Cells.Item(105, 5).Value = 3
$wsBR.Cells.Item(105, 6).Value = "Occupied"

$wsBR.Cells.Item(106, 1).NumberFormat = "@"
$wsBR.Cells.Item(106, 1).Value = "2026-02-01"
$wsBR.Cells.Item(106, 2).Value = "21:07:15"
$wsBR.Cells.Item(106, 3).Value = "21:00"
$wsBR.Cells.Item(106, 4).Value = "Bedroom"
$wsBR.Cells.Item(106, 5).Value = 2
$wsBR.Cells.Item(106, 6).Value = "Occupied"

$wsBR.Cells.Item(107, 1).NumberFormat = "@"
$wsBR.Cells.Item(107, 1).Value = "2026-02-01"
$wsBR.Cells.Item(107, 2).Value = "21:07:21"
$wsBR.Cells.Item(107, 3).Value = "21:00"
$wsBR.Cells.Item(107, 4).Value = "Bedroom"
$wsBR.Cells.Item(107, 5).Value = 1
$wsBR.Cells.Item(107, 6).Value = "Occupied"

$wsBR.Cells.Item(108, 1).NumberFormat = "@"
$wsBR.Cells.Item(108, 1).Value = "2026-02-01"
$wsBR.Cells.Item(108, 2).Value = "21:07:28"
$wsBR.Cells.Item(108, 3).Value = "21:00"
$wsBR.Cells.Item(108, 4).Value = "Bedroom"
$wsBR.Cells.Item(108, 5).Value = 2
$wsBR.Cells.Item(108, 6).Value = "Occupied"

$wsBR.Cells.Item(109, 1).NumberFormat = "@"
$wsBR.Cells.Item(109, 1).Value = "2026-02-01"
$wsBR.Cells.Item(109, 2).Value = "21:07:31"
$wsBR.Cells.Item(109, 3).Value = "21:00"
$wsBR.Cells.Item(109, 4).Value = "Bedroom"
$wsBR.Cells.Item(109, 5).Value = 1
$wsBR.Cells.Item(109, 6).Value = "Occupied"

$wsBR.Cells.Item(110, 1).NumberFormat = "@"
$wsBR.Cells.Item(110, 1).Value = "2026-02-01"
$wsBR.Cells.Item(110, 2).Value = "21:07:37"
$wsBR.Cells.Item(110, 3).Value = "21:00"
$wsBR.Cells.Item(110, 4).Value = "Bedroom"
$wsBR.Cells.Item(110, 5).Value = 2
$wsBR.Cells.Item(110, 6).Value = "Occupied"

$wsBR.Cells.Item(111, 1).NumberFormat = "@"
$wsBR.Cells.Item(111, 1).Value = "2026-02-01"
$wsBR.Cells.Item(111, 2).Value = "21:07:41"
$wsBR.Cells.Item(111, 3).Value = "21:00"
$wsBR.Cells.Item(111, 4).Value = "Bedroom"
$wsBR.Cells.Item(111, 5).Value = 1
$wsBR.Cells.Item(111, 6).Value = "Occupied"

# --- mmWave(HR): append rows 97-112 ---
$wsHR = $wb.Worksheets.Item("mmWave(HR)")
$wsHR.Cells.Item(97, 1).NumberFormat = "@"
$wsHR.Cells.Item(97, 1).Value = "2026-02-01"
$wsHR.Cells.Item(97, 2).Value = "21:06:43"
$wsHR.Cells.Item(97, 3).Value = "21:00"
$wsHR.Cells.Item(97, 4).Value = "Bedroom"
$wsHR.Cells.Item(97, 5).Value = 50
$wsHR.Cells.Item(97, 6).Value = "Occupied"

$wsHR.Cells.Item(98, 1).NumberFormat = "@"
$wsHR.Cells.Item(98, 1).Value = "2026-02-01"
$wsHR.Cells.Item(98, 2).Value = "21:06:47"
$wsHR.Cells.Item(98, 3).Value = "21:00"
$wsHR.Cells.Item(98, 4).Value = "Bedroom"
$wsHR.Cells.Item(98, 5).Value = 49
$wsHR.Cells.Item(98, 6).Value = "Occupied"

$wsHR.Cells.Item(99, 1).NumberFormat = "@"
$wsHR.Cells.Item(99, 1).Value = "2026-02-01"
$wsHR.Cells.Item(99, 2).Value = "21:06:55"
$wsHR.Cells.Item(99, 3).Value = "21:00"
$wsHR.Cells.Item(99, 4).Value = "Bedroom"
$wsHR.Cells.Item(99, 5).Value = 50
$wsHR.Cells.Item(99, 6).Value = "Occupied"

$wsHR.Cells.Item(100, 1).NumberFormat = "@"
$wsHR.Cells.Item(100, 1).Value = "2026-02-01"
$wsHR.Cells.Item(100, 2).Value = "21:06:59"
$wsHR.Cells.Item(100, 3).Value = "21:00"
$wsHR.Cells.Item(100, 4).Value = "Bedroom"
$wsHR.Cells.Item(100, 5).Value = 51
$wsHR.Cells.Item(100, 6).Value = "Occupied"

$wsHR.Cells.Item(101, 1).NumberFormat = "@"
$wsHR.Cells.Item(101, 1).Value = "2026-02-01"
$wsHR.Cells.Item(101, 2).Value = "21:07:01"
$wsHR.Cells.Item(101, 3).Value = "21:00"
$wsHR.Cells.Item(101, 4).Value = "Bedroom"
$wsHR.Cells.Item(101, 5).Value = 50
$wsHR.Cells.Item(101, 6).Value = "Occupied"

$wsHR.Cells.Item(102, 1).NumberFormat = "@"
$wsHR.Cells.Item(102, 1).Value = "2026-02-01"
$wsHR.Cells.Item(102, 2).Value = "21:07:04"
$wsHR.Cells.Item(102, 3).Value = "21:00"
$wsHR.Cells.Item(102, 4).Value = "Bedroom"
$wsHR.Cells.Item(102, 5).Value = 49
$wsHR.Cells.Item(102, 6).Value = "Occupied"

$wsHR.Cells.Item(103, 1).NumberFormat = "@"
$wsHR.Cells.Item(103, 1).Value = "2026-02-01"
$wsHR.Cells.Item(103, 2).Value = "21:07:05"
$wsHR.Cells.Item(103, 3).Value = "21:00"
$wsHR.Cells.Item(103, 4).Value = "Bedroom"
$wsHR.Cells.Item(103, 5).Value = 50
$wsHR.Cells.Item(103, 6).Value = "Occupied"

$wsHR.Cells.Item(104, 1).NumberFormat = "@"
$wsHR.Cells.Item(104, 1).Value = "2026-02-01"
$wsHR.Cells.Item(104, 2).Value = "21:07:08"
$wsHR.Cells.Item(104, 3).Value = "21:00"
$wsHR.Cells.Item(104, 4).Value = "Bedroom"
$wsHR.Cells.Item(104, 5).Value = 49
$wsHR.Cells.Item(104, 6).Value = "Occupied"

$wsHR.Cells.Item(105, 1).NumberFormat = "@"
$wsHR.Cells.Item(105, 1).Value = "2026-02-01"
$wsHR.Cells.Item(105, 2).Value = "21:07:13"
$wsHR.Cells.Item(105, 3).Value = "21:00"
$wsHR.Cells.Item(105, 4).Value = "Bedroom"
$wsHR.Cells.Item(105, 5).Value = 50
$wsHR.Cells.Item(105, 6).Value = "Occupied"

$wsHR.Cells.Item(106, 1).NumberFormat = "@"
$wsHR.Cells.Item(106, 1).Value = "2026-02-01"
$wsHR.Cells.Item(106, 2).Value = "21:07:14"
$wsHR.Cells.Item(106, 3).Value = "21:00"
$wsHR.Cells.Item(106, 4).Value = "Bedroom"
$wsHR.Cells.Item(106, 5).Value = 51
$wsHR.Cells.Item(106, 6).Value = "Occupied"

$wsHR.Cells.Item(107, 1).NumberFormat = "@"
$wsHR.Cells.Item(107, 1).Value = "2026-02-01"
$wsHR.Cells.Item(107, 2).Value = "21:07:15"
$wsHR.Cells.Item(107, 3).Value = "21:00"
$wsHR.Cells.Item(107, 4).Value = "Bedroom"
$wsHR.Cells.Item(107, 5).Value = 50
$wsHR.Cells.Item(107, 6).Value = "Occupied"

$wsHR.Cells.Item(108, 1).NumberFormat = "@"
$wsHR.Cells.Item(108, 1).Value = "2026-02-01"
$wsHR.Cells.Item(108, 2).Value = "21:07:20"
$wsHR.Cells.Item(108, 3).Value = "21:00"
$wsHR.Cells.Item(108, 4).Value = "Bedroom"
$wsHR.Cells.Item(108, 5).Value = 49
$wsHR.Cells.Item(108, 6).Value = "Occupied"

$wsHR.Cells.Item(109, 1).NumberFormat = "@"
$wsHR.Cells.Item(109, 1).Value = "2026-02-01"
$wsHR.Cells.Item(109, 2).Value = "21:07:28"
$wsHR.Cells.Item(109, 3).Value = "21:00"
$wsHR.Cells.Item(109, 4).Value = "Bedroom"
$wsHR.Cells.Item(109, 5).Value = 50
$wsHR.Cells.Item(109, 6).Value = "Occupied"

$wsHR.Cells.Item(110, 1).NumberFormat = "@"
$wsHR.Cells.Item(110, 1).Value = "2026-02-01"
$wsHR.Cells.Item(110, 2).Value = "21:07:31"
$wsHR.Cells.Item(110, 3).Value = "21:00"
$wsHR.Cells.Item(110, 4).Value = "Bedroom"
$wsHR.Cells.Item(110, 5).Value = 49
$wsHR.Cells.Item(110, 6).Value = "Occupied"

$wsHR.Cells.Item(111, 1).NumberFormat = "@"
$wsHR.Cells.Item(111, 1).Value = "2026-02-01"
$wsHR.Cells.Item(111, 2).Value = "21:07:37"
$wsHR.Cells.Item(111, 3).Value = "21:00"
$wsHR.Cells.Item(111, 4).Value = "Bedroom"
$wsHR.Cells.Item(111, 5).Value = 50
$wsHR.Cells.Item(111, 6).Value = "Occupied"

$wsHR.Cells.Item(112, 1).NumberFormat = "@"
$wsHR.Cells.Item(112, 1).Value = "2026-02-01"
$wsHR.Cells.Item(112, 2).Value = "21:07:40"
$wsHR.Cells.Item(112, 3).Value = "21:00"
$wsHR.Cells.Item(112, 4).Value = "Bedroom"
$wsHR.Cells.Item(112, 5).Value = 49
$wsHR.Cells.Item(112, 6).Value = "Occupied"

